# Merge the duplicate "PARTIDO ..." rows (29-35) back into their matching
# short-named rows (7, 8, 14, 20, 21, 24, 26), then delete the now-redundant
# rows 29-35 entirely, shrinking the used range from A1:F35 to A1:F28.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7  (DE LA GENTE)            <- old row 29 (PARTIDO DE LA GENTE)
$ws.Range("C7").Value = 0.5
$ws.Range("D7").Value = 0.4459459459459459
$ws.Range("F7").Value = 0.3333333333333333

# Row 8  (DEMOCRATA CRISTIANO)    <- old row 30 (PARTIDO DEMOCRATA CRISTIANO)
$ws.Range("C8").Value = 0.4545454545454545
$ws.Range("D8").Value = 0.4523809523809524
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0.125

# Row 14 (LIBERAL DE CHILE)       <- old row 31 (PARTIDO LIBERAL DE CHILE)
$ws.Range("C14").Value = 0.5
$ws.Range("D14").Value = 0.4117647058823529
$ws.Range("F14").Value = 0

# Row 20 (POR LA DEMOCRACIA)      <- old row 32 (PARTIDO POR LA DEMOCRACIA)
$ws.Range("C20").Value = 0.6
$ws.Range("D20").Value = 0.4545454545454545
$ws.Range("E20").Value = 0.5
$ws.Range("F20").Value = 0.4285714285714285

# Row 21 (RADICAL DE CHILE)       <- old row 33 (PARTIDO RADICAL DE CHILE)
$ws.Range("C21").Value = 0.3333333333333333
$ws.Range("D21").Value = 0.4642857142857143
$ws.Range("F21").Value = 0

# Row 24 (REPUBLICANO DE CHILE)   <- old row 34 (PARTIDO REPUBLICANO DE CHILE)
$ws.Range("C24").Value = 0.5
$ws.Range("D24").Value = 0.4142857142857143
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 0.1428571428571428

# Row 26 (SOCIALISTA DE CHILE)    <- old row 35 (PARTIDO SOCIALISTA DE CHILE)
$ws.Range("C26").Value = 0.4444444444444444
$ws.Range("D26").Value = 0.4418604651162791
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 0.3076923076923077

# The data that used to live in rows 29-35 is now merged above; delete the
# now-redundant trailing rows so the sheet shrinks from A1:F35 to A1:F28.
$ws.Range("A29:F35").EntireRow.Delete()
